$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 11-43 first (they get folded into the consolidated rows 2-10)
$ws.Range("A11:A43").EntireRow.Delete() | Out-Null

# Now set the consolidated text for rows 2-10
$ws.Range("A2").Value = "('Blightning', ['{1}{B}{R}', 'Sorcery', 'Blightning deals 3 damage to target player or planeswalker. That player or that planeswalker" + [char]8217 + "s controller discards two cards.'])"
$ws.Range("A3").Value = "('Cryptic Command', ['{1}{U}{U}{U}', 'Instant', 'Choose two " + [char]8212 + "', '" + [char]8226 + " Counter target spell.', '" + [char]8226 + " Return target permanent to its owner" + [char]8217 + "s hand.', '" + [char]8226 + " Tap all creatures your opponents control.', '" + [char]8226 + " Draw a card.'])"
$ws.Range("A4").Value = "('Flame Javelin', ['{2/R}{2/R}{2/R}', 'Instant', '({2/R} can be paid with any two mana or with {R}. This card" + [char]8217 + "s converted mana cost is 6.)', 'Flame Javelin deals 4 damage to any target.'])"
$ws.Range("A5").Value = "('Nameless Inversion', ['{1}{B}', 'Tribal Instant " + [char]8212 + " Shapeshifter', 'Changeling (This card is every creature type.)', 'Target creature gets +3/-3 and loses all creature types until end of turn.'])"
$ws.Range("A6").Value = "('Negate', ['{1}{U}', 'Instant', 'Counter target noncreature spell.'])"
$ws.Range("A7").Value = "('Rampant Growth', ['{1}{G}', 'Sorcery', 'Search your library for a basic land card and put that card onto the battlefield tapped. Then shuffle your library.'])"
$ws.Range("A8").Value = "('Remove Soul', ['{1}{U}', 'Instant', 'Counter target creature spell.'])"
$ws.Range("A9").Value = "('Terminate', ['{B}{R}', 'Instant', 'Destroy target creature. It can" + [char]8217 + "t be regenerated.'])"
$ws.Range("A10").Value = "('Unmake', ['{W/B}{W/B}{W/B}', 'Instant', 'Exile target creature.'])"
